$d = $word.ActiveDocument

# The trailing "_GoBack" bookmark currently sits right after "Mich wunderts
# nicht." (inside that paragraph, just before its paragraph mark). It needs
# to end up after the very last line of the new text instead. Remove it now;
# it will be re-created (with the same id) in the right spot further down.
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
}

# The document currently ends with an empty paragraph (it only carries a
# de-DE language tag). Replace that paragraph outright with the three new
# paragraphs of poetry, preserving/assigning the same run/paragraph
# formatting pattern already used throughout the document, and re-attach
# the _GoBack bookmark (collapsed) right after the final run of text.
$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$target = $lastParagraph.Range

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$newXml = "<w:p $w><w:pPr><w:rPr><w:rtl/><w:lang w:val='de-DE'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='de-DE'/></w:rPr><w:t>Das Gesetz hat noch keinen großen Mann gebildet,</w:t></w:r></w:p>" + `
          "<w:p $w><w:pPr><w:rPr><w:lang w:val='de-DE'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='de-DE'/></w:rPr><w:t>die Freiheit brütet Kolosse und Extremitäten aus</w:t></w:r></w:p>" + `
          "<w:p $w><w:pPr><w:rPr><w:lang w:val='de-DE'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='de-DE'/></w:rPr><w:t>Ich bin der Glücklichste unter der Sonne</w:t></w:r>" + `
          "<w:bookmarkStart $w w:id='0' w:name='_GoBack'/><w:bookmarkEnd $w w:id='0'/></w:p>"

$target.InsertXML($newXml)
